# Add a new "Trajectory Ordering" section at the bottom of Sheet1,
# mirroring the style of the existing section headers (e.g. "Debug").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Section header row (row 34) - same look as A2/A9/A22/A31 ("Debug", etc.)
$ws.Range("A34").Value = "Trajectory Ordering"
$ws.Range("A34").Font.Underline = $true

# Parameter row (row 35)
$ws.Range("A35").Value = "Contour First?"
$ws.Range("B35").Value = "Yes"

# Scroll the sheet view so row 17 is at the top, as in the saved workbook.
$excel.ActiveWindow.ScrollRow = 17
